# Update gh-pages to output generated at 456a3b4
# Applies updated "want to go" counts (F) and minimum ticket price (G)
# figures to the "展览" (sheet 1) and "全部类型" (sheet 4) worksheets.

$wb = $excel.ActiveWorkbook

function Update-ConListSheet($ws, $Row6, $Row7, $Row11, $Row15, $Row17, $Row18, $Row19, $Row20, $Row23, $Row26, $Row31, $Row36, $Row39, $Row42) {
    $ws.Range("F4").Value = 14006
    $ws.Range("F5").Value = 242

    $ws.Range("F$Row6").Value = 1806
    $ws.Range("G$Row6").Value = "不可售"

    $ws.Range("F$Row7").Value = 172
    $ws.Range("G$Row7").Value = "不可售"

    $ws.Range("F8").Value = 149
    $ws.Range("G8").Value = 49.9

    $ws.Range("F$Row11").Value = 560
    $ws.Range("F$Row15").Value = 14177
    $ws.Range("F$Row17").Value = 638
    $ws.Range("F$Row18").Value = 15038
    $ws.Range("F$Row19").Value = 19
    $ws.Range("F$Row20").Value = 8392
    $ws.Range("F$Row23").Value = 40
    $ws.Range("F$Row26").Value = 170
    $ws.Range("F$Row31").Value = 1047
    $ws.Range("F$Row36").Value = 13
    $ws.Range("F$Row39").Value = 233
    $ws.Range("F$Row42").Value = 5166
}

# Sheet "展览" (rows 6,7,11,15,17,18,19,20,23,26,31,36,39,42 line up directly)
$wsExpo = $wb.Worksheets.Item("展览")
Update-ConListSheet $wsExpo 6 7 11 15 17 18 19 20 23 26 31 36 39 42

# Sheet "全部类型" mirrors the same records, but rows 36/39/42 are shifted
# down by two (38/41/44) because of extra rows present only in this sheet.
$wsAll = $wb.Worksheets.Item("全部类型")
Update-ConListSheet $wsAll 6 7 11 15 17 18 19 20 23 26 31 38 41 44
